$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing data row is 289 (dates run through 2025-07-07, serial 45845).
# Append three more rows (290-292) that repeat the same data but with the
# date serial incremented by one day each time (45846, 45847, 45848),
# matching the style/number-format used by the existing date column.

$srcRow = 289
$newRows = @(290, 291, 292)

foreach ($r in $newRows) {
    # Copy the whole source row's formatting (styles) into the new row
    # without touching values yet, so we reuse the existing style records
    # instead of creating new duplicate ones.
    $ws.Range("A$srcRow`:J$srcRow").Copy()
    $ws.Range("A$r`:J$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

$offset = 1
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2 + $offset
    for ($col = 2; $col -le 10; $col++) {
        $ws.Cells.Item($r, $col).Value2 = $ws.Cells.Item($srcRow, $col).Value2
    }
    $offset = $offset + 1
}
